$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the leading zero: without forcing text format Excel would
# auto-convert "001" to the number 1.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

$ws.Range("M2").Value = "2020-12-17 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = -250642402.73
$ws.Range("P2").Value = -92.76719467380001
$ws.Range("Q2").Value = 4270700072.05
$ws.Range("R2").Value = 1580.6617741537
$ws.Range("S2").Value = 204432168.31
$ws.Range("T2").Value = 75.6639680622
$ws.Range("U2").Value = -604662318.52
$ws.Range("V2").Value = -223.7962387972
$ws.Range("W2").Value = 1542591.36
$ws.Range("X2").Value = 0.5709403973
$ws.Range("Y2").Value = 514082712.59
$ws.Range("Z2").Value = 190.2711215574
$ws.Range("AA2").Value = 1123180658.18
$ws.Range("AB2").Value = 415.7090645332
$ws.Range("AC2").Value = 270184307.73
$ws.Range("AD2").Value = 87.681237337
